# "Hjemme passive tweaks lichtwark deleted values"
# Overwrite the first four data columns (B:E) of the two data rows with the
# "lichtwark"-adjusted trial-max values, and update the header row (B1:E1)
# to match the corresponding trial lengths. Also restyle the live selection
# down to the edited block (B1:E3) instead of the full original range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - trial length headers
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - CON trial max values
$ws.Range("B2").Value = 265.6287409710601
$ws.Range("C2").Value = 290.73128431545058
$ws.Range("D2").Value = 262.21989162037539
$ws.Range("E2").Value = 297.51363840877923

# Row 3 - STR trial max values
$ws.Range("B3").Value = 261.16854873030132
$ws.Range("C3").Value = 295.99799930273218
$ws.Range("D3").Value = 261.79447028965973
$ws.Range("E3").Value = 303.44048945715053

# Shrink the saved selection to the edited block
$ws.Range("B1:E3").Select()
